$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 321
$ws.Range("F3").Value = 1120
$ws.Range("F4").Value = 1241
$ws.Range("F5").Value = 1121
$ws.Range("F6").Value = 3381
$ws.Range("F10").Value = 766
$ws.Range("F11").Value = 588
$ws.Range("F12").Value = 54
$ws.Range("F14").Value = 649
$ws.Range("F15").Value = 1769
$ws.Range("F17").Value = 359
$ws.Range("F18").Value = 34
$ws.Range("F19").Value = 55
$ws.Range("F20").Value = 658
$ws.Range("F21").Value = 412
$ws.Range("F22").Value = 740
$ws.Range("F23").Value = 79164
$ws.Range("F24").Value = 79164
$ws.Range("F27").Value = 33673
$ws.Range("F28").Value = 33673
$ws.Range("F29").Value = 503
$ws.Range("F33").Value = 34
$ws.Range("F34").Value = 969
$ws.Range("F35").Value = 295
$ws.Range("F37").Value = 585
$ws.Range("F38").Value = 1066
$ws.Range("F40").Value = 5462
$ws.Range("F41").Value = 766
$ws.Range("F42").Value = 449
$ws.Range("F46").Value = 385
$ws.Range("F50").Value = 45

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 19
$ws.Range("F14").Value = 1770
$ws.Range("F15").Value = 25
$ws.Range("F18").Value = 412
$ws.Range("F20").Value = 71
$ws.Range("F21").Value = 74
$ws.Range("F23").Value = 518
$ws.Range("F24").Value = 518
$ws.Range("F26").Value = 769
$ws.Range("F47").Value = 68
$ws.Range("F48").Value = 824
$ws.Range("F49").Value = 139
$ws.Range("F51").Value = 63

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 727
$ws.Range("F5").Value = 560
$ws.Range("F6").Value = 596
$ws.Range("F7").Value = 98

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 727
$ws.Range("F3").Value = 560
$ws.Range("F5").Value = 1241
$ws.Range("F7").Value = 1121
$ws.Range("F8").Value = 3381
$ws.Range("F11").Value = 766
$ws.Range("F12").Value = 596
$ws.Range("F14").Value = 588
$ws.Range("F15").Value = 54
$ws.Range("F16").Value = 649
$ws.Range("F17").Value = 98
$ws.Range("F18").Value = 1769
$ws.Range("F19").Value = 25
$ws.Range("F23").Value = 34
$ws.Range("F24").Value = 55
$ws.Range("F25").Value = 658
$ws.Range("F26").Value = 412
$ws.Range("F27").Value = 412
$ws.Range("F28").Value = 71
$ws.Range("F29").Value = 79164
$ws.Range("F30").Value = 74
$ws.Range("F31").Value = 33673
$ws.Range("F32").Value = 503
$ws.Range("F35").Value = 518
$ws.Range("F36").Value = 34
$ws.Range("F38").Value = 969
$ws.Range("F41").Value = 295
$ws.Range("F44").Value = 585
$ws.Range("F45").Value = 1066
$ws.Range("F47").Value = 766
$ws.Range("F49").Value = 449
$ws.Range("F53").Value = 68
$ws.Range("F54").Value = 139
